$d = $word.ActiveDocument

# 1. Gabinete paragraph: add conditional S/N block before "de su interior se extrae:"
$d.Content.Find.Execute(
    "Un (01) gabinete, color XX, marca {marca}, modelo {modelo}, de su interior se extrae:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Un (01) gabinete, color XX, marca {marca}, modelo {modelo},{#serialNumber==“”} Sin S/N visible{/}{#serialNumber!=“”} con S/N: {serialNumber}{/}, de su interior se extrae:",
    2
)

# 2. Disk paragraphs (two identical occurrences): replace Tableau duplicator wording
#    with generic forensic software wording using {herramientaSoftDisco}.
$d.Content.Find.Execute(
    "se procedió a conectarlo a un duplicador forense marca Tableau, a los fines de realizar una adquisición forense de su contenido{#estadoDisco==”completo”} a través del software Tableau Imager,",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "se procedió a realizar una tarea técnica utilizando el software forense aportado por la empresa {herramientaSoftDisco}, a los fines de realizar una adquisición forense de su contenido{#estadoDisco==”completo”}",
    2
)

Write-Output "done"
